$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Preço Atual" (current price) values for rows 2-11 (column C)
$newPrices = @{
    2  = 84.66
    3  = 158.83
    4  = 277.8
    5  = 386.34
    6  = 482.28
    7  = 134.68
    8  = 1106.86
    9  = 333.47
    10 = 1546.16
    11 = 9.01
}

foreach ($row in $newPrices.Keys) {
    $idealPrice = $ws.Cells.Item($row, 2).Value2
    $currentPrice = $newPrices[$row]

    $ws.Cells.Item($row, 3).Value = $currentPrice
    $ws.Cells.Item($row, 4).Value = [bool]($currentPrice -lt $idealPrice)
}
